# Updated symbol list — refreshes Price (D) and Volume(1h) (E) columns for
# this row of coins to the latest scraped quotes.
#
# Target cells are stored as plain text (inlineStr, General format) in the
# workbook, e.g. "327.12" / "0.16%", not as numbers/percentages. Writing a
# numeric-looking string straight into .Value lets Excel's type-inference
# convert it to a real number (and restyle it). To keep these as literal
# text we flip NumberFormat to "@" (Text) right before the assignment, then
# restore the cell style to "Normal" afterwards so no stray text-format
# style lingers on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2" = "327.12"
    "E2" = "0.16%"
    "D3" = "44.02"
    "E3" = "-2.82%"
    "D4" = "5.512"
    "E4" = "-0.98%"
    "E5" = "-0.85%"
    "D8" = "0.9480"
    "E8" = "-0.17%"
    "D9" = "0.1146"
    "E9" = "-1.31%"
    "D10" = "0.1838"
    "E10" = "-3.02%"
    "D11" = "12.19"
    "E11" = "42.60%"
    "D12" = "0.09755"
    "E12" = "-3.04%"
    "D13" = "0.04611"
    "E13" = "10.53%"
    "E14" = "0.20%"
    "D15" = "0.001267"
    "E15" = "-0.68%"
    "D16" = "0.04079"
    "E16" = "-4.58%"
    "D17" = "0.005721"
    "E17" = "-5.85%"
    "E18" = "-6.96%"
    "D19" = "4.288"
    "E19" = "-1.18%"
    "E20" = "-0.20%"
    "D21" = "0.1406"
    "E21" = "2.28%"
    "D23" = "0.001244"
    "E23" = "0.49%"
    "D24" = "0.004320"
    "E24" = "-6.04%"
    "D25" = "0.0001188"
    "E25" = "-3.70%"
    "D26" = "0.0003742"
    "E26" = "-6.57%"
    "D38" = "0.02565"
    "E38" = "-4.27%"
    "D39" = "0.05562"
    "E39" = "0.12%"
    "D40" = "0.007544"
    "E40" = "-1.86%"
    "E41" = "0.08%"
    "D42" = "0.007601"
    "E42" = "-32.94%"
    "D43" = "0.002012"
    "E43" = "-3.04%"
    "D44" = "0.008513"
    "E44" = "-2.19%"
    "D45" = "0.00007120"
    "E45" = "-0.01%"
    "D46" = "0.00000000750"
    "E46" = "-0.34%"
    "D47" = "0.003528"
    "E47" = "54.80%"
    "D48" = "0.001845"
    "E48" = "-46.46%"
    "D49" = "0.00002099"
    "E49" = "-0.34%"
    "D50" = "0.0001999"
    "E50" = "-0.34%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Style = "Normal"
}
